$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for several rows - repulled data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = 1
